# Update course schedule: updating schedule; post recording
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Week 7 row (row 8): split "Uncertainty & intro to tables and fonts" into
#     a standalone "Uncertainty" theme, trim its description, and add the
#     newly-posted lecture recording link.
$ws.Range("C8").Value = 'Uncertainty'
$ws.Range("E8").Value = 'Common methods for visualizing uncertainty (and their implementation w/{ggplot2}). Framing uncertainty as relative frequencies. Non-standard methods for visualizing standard errors, boostrapping, and a brief foray into hypothetical outcomes plots. We''ll also discuss changing fonts, both within websites/applications, as well as with {ggplot2}.'
$ws.Range("M8").Value = 'https://youtu.be/uYj04BODzIc'

# --- Week 8 row (row 9): rename theme/description to fold in the distill +
#     fonts content that used to live on week 7; the HW2 due date moves off
#     this row (draft only remains due here).
$ws.Range("C9").Value = 'Websites, flex dashbaords, fonts, and some customization with CSS'
$ws.Range("E9").Value = 'Websites with [{distill}](https://rstudio.github.io/distill/), which help you create relatively simple yet customizable blogs, optimized for scientific communication. Also building (static) data dashboards with the [{flexdashboard}](https://rmarkdown.rstudio.com/flexdashboard/) package. Finally, we will discuss customization with CSS, and changing the fonts in both web-based documents, as well as ggplot2 plots.'
$ws.Range("H9").Value = 'assignments/#draft'
$ws.Range("I9").Value = 'Draft'

# --- Week 9 row (row 10): theme/description now also covers tables content
#     that moved off week 7; HW2 becomes due on this row alongside the peer
#     review.
$ws.Range("C10").Value = 'Tables & intro to Geographic data'
$ws.Range("E10").Value = 'We will focus primarily on two packages for creating tables: [{gt}](https://gt.rstudio.com) for static tables, and [{reactable}](https://glin.github.io/reactable/index.html) for interactive tables. We will also discuss the differences between vector and raster data, producing basic maps, getting data for producing different types of maps, and understandin the basics of the R geospatial ecosystem (which is consistently and rapidly evolving).'
$ws.Range("H10").Value = 'assignments/#peer-review; homework-2'
$ws.Range("I10").Value = 'PR; HW2'

# Row heights grew to fit the new wrapped text.
$ws.Rows.Item(9).RowHeight = 154
$ws.Rows.Item(10).RowHeight = 136

# Move the current selection to reflect where the author was last working.
$ws.Range("C10").Select()
